$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.906.10'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '1.777.20'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'315.41"
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = "'0.5383"
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('D8').Value = "'0.3754"
$ws.Range('E8').Value = '  -2.85%  '
$ws.Range('D9').Value = "'0.07432"
$ws.Range('E9').Value = '  -1.96%  '
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').Value = "'1.091"
$ws.Range('E11').Value = '  -2.60%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = "'20.39"
$ws.Range('E13').Value = '  -3.30%  '
$ws.Range('D14').Value = "'6.064"
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('D16').Value = '1.771.18'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').Value = "'87.87"
$ws.Range('E17').Value = '  -4.69%  '
$ws.Range('D18').Value = "'0.00001052"
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').Value = "'0.06410"
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').Value = "'1.000"
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = "'5.869"
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('D23').Value = '27.936.78'
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('D24').Value = "'11.13"
$ws.Range('E24').Value = '  -2.54%  '
$ws.Range('E25').Value = '  -2.31%  '
$ws.Range('D26').Value = "'155.61"
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('D27').Value = "'20.20"
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').Value = '1.973.44'
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').Value = "'2.270"
$ws.Range('E29').Value = '  -5.20%  '
$ws.Range('E30').Value = '  -3.25%  '
$ws.Range('D31').Value = "'1.108"
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('D32').Value = "'0.1051"
$ws.Range('E32').Value = '  +2.94%  '
$ws.Range('D33').Value = "'3.642"
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D34').Value = "'5.506"
$ws.Range('E34').Value = '  -3.94%  '
$ws.Range('D35').Value = "'0.2240"
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('D36').Value = "'0.06359"
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('D37').Value = "'0.02259"
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').Value = "'4.953"
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D39').Value = "'8.372"
$ws.Range('E39').Value = '  -5.74%  '
$ws.Range('D40').Value = "'0.6111"
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = "'10.99"
$ws.Range('E41').Value = '  -5.31%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.175"
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').Value = "'1.426"
$ws.Range('E43').Value = '  +3.27%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = "'13.30"
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = "'3.649"
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = "'0.5733"
$ws.Range('E47').Value = '  -3.95%  '
$ws.Range('D48').Value = "'125.78"
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').Value = "'1.181"
$ws.Range('E49').Value = '  +3.06%  '
$ws.Range('D50').Value = "'1.918"
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('E51').Value = '  -1.80%  '
